# "Added New Mac-Address and Document Types"
# Appends 5 new device rows (device_id 3000176-3000180) for
# regcntr_id 10002 / machine_id 10032, mirroring the existing rows'
# pattern (lang_code "eng", is_active TRUE, cr_by "superadmin",
# cr_dtimes/eff_dtimes "now()"), and switches the workbook to manual
# calculation mode.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Workbook now recalculates manually (calcPr calcMode="manual")
$excel.Calculation = -4135   # xlCalculationManual

$startRow  = 157
$regcntrId = 10002
$machineId = 10032
$deviceIds = @(3000176, 3000177, 3000178, 3000179, 3000180)

for ($i = 0; $i -lt $deviceIds.Length; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value = $regcntrId
    $ws.Cells.Item($row, 2).Value = $machineId
    $ws.Cells.Item($row, 3).Value = $deviceIds[$i]
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin"
    $ws.Cells.Item($row, 7).Value = "now()"
    $ws.Cells.Item($row, 8).Value = "now()"
}

# Reflect the scrolled/selected state left behind after the edit
$excel.ActiveWindow.ScrollRow = 151
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D157").Select() | Out-Null
